$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 18.741769
$ws.Range("H2").Value = 56.225307
$ws.Range("I2").Value = 0.2218531826860132
$ws.Range("J2").Value = 0.2218531826860132
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.521994666666667
$ws.Range("N2").Value = 7.565983999999999
$ws.Range("O2").Value = 0.01218715015226367
$ws.Range("P2").Value = 0.01218715015226367
$ws.Range("Q2").Value = 47.26664146189867
$ws.Range("R2").Value = 425.399773157088
$ws.Range("S2").Value = 0.002703758049152026
$ws.Range("T2").Value = 0.002703758049152026
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 18.741769
$ws.Range("H3").Value = 56.225307
$ws.Range("I3").Value = 0.2218531826860132
$ws.Range("J3").Value = 0.2218531826860132
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 157.1645253333333
$ws.Range("N3").Value = 471.493576
$ws.Range("O3").Value = 0.7594733225102963
$ws.Range("P3").Value = 0.7594733225102964
$ws.Range("Q3").Value = 2945.541228791982
$ws.Range("R3").Value = 26509.87105912783
$ws.Range("S3").Value = 0.1684915737640302
$ws.Range("T3").Value = 0.1684915737640302
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 18.741769
$ws.Range("H4").Value = 56.225307
$ws.Range("I4").Value = 0.2218531826860132
$ws.Range("J4").Value = 0.2218531826860132
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 47.252316
$ws.Range("N4").Value = 141.756948
$ws.Range("O4").Value = 0.2283395273374399
$ws.Range("P4").Value = 0.2283395273374399
$ws.Range("Q4").Value = 885.5919911870041
$ws.Range("R4").Value = 7970.327920683036
$ws.Range("S4").Value = 0.05065785087283097
$ws.Range("T4").Value = 0.05065785087283097
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 53.77230066666667
$ws.Range("H5").Value = 161.316902
$ws.Range("I5").Value = 0.6365224138259964
$ws.Range("J5").Value = 0.6365224138259964
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.521994666666667
$ws.Range("N5").Value = 7.565983999999999
$ws.Range("O5").Value = 0.01218715015226367
$ws.Range("P5").Value = 0.01218715015226367
$ws.Range("Q5").Value = 135.6134554957298
$ws.Range("R5").Value = 1220.521099461568
$ws.Range("S5").Value = 0.007757394232578731
$ws.Range("T5").Value = 0.007757394232578731
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 53.77230066666667
$ws.Range("H6").Value = 161.316902
$ws.Range("I6").Value = 0.6365224138259964
$ws.Range("J6").Value = 0.6365224138259964
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 157.1645253333333
$ws.Range("N6").Value = 471.493576
$ws.Range("O6").Value = 0.7594733225102963
$ws.Range("P6").Value = 0.7594733225102964
$ws.Range("Q6").Value = 8451.098110357951
$ws.Range("R6").Value = 76059.88299322156
$ws.Range("S6").Value = 0.4834217924807033
$ws.Range("T6").Value = 0.4834217924807033
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 53.77230066666667
$ws.Range("H7").Value = 161.316902
$ws.Range("I7").Value = 0.6365224138259964
$ws.Range("J7").Value = 0.6365224138259964
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 47.252316
$ws.Range("N7").Value = 141.756948
$ws.Range("O7").Value = 0.2283395273374399
$ws.Range("P7").Value = 0.2283395273374399
$ws.Range("Q7").Value = 2540.865743148344
$ws.Range("R7").Value = 22867.7916883351
$ws.Range("S7").Value = 0.1453432271127144
$ws.Range("T7").Value = 0.1453432271127144
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 11.964182
$ws.Range("H8").Value = 35.892546
$ws.Range("I8").Value = 0.1416244034879904
$ws.Range("J8").Value = 0.1416244034879904
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.521994666666667
$ws.Range("N8").Value = 7.565983999999999
$ws.Range("O8").Value = 0.01218715015226367
$ws.Range("P8").Value = 0.01218715015226367
$ws.Range("Q8").Value = 30.17360319502933
$ws.Range("R8").Value = 271.562428755264
$ws.Range("S8").Value = 0.001725997870532914
$ws.Range("T8").Value = 0.001725997870532914
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 11.964182
$ws.Range("H9").Value = 35.892546
$ws.Range("I9").Value = 0.1416244034879904
$ws.Range("J9").Value = 0.1416244034879904
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 157.1645253333333
$ws.Range("N9").Value = 471.493576
$ws.Range("O9").Value = 0.7594733225102963
$ws.Range("P9").Value = 0.7594733225102964
$ws.Range("Q9").Value = 1880.344985031611
$ws.Range("R9").Value = 16923.10486528449
$ws.Range("S9").Value = 0.1075599562655629
$ws.Range("T9").Value = 0.1075599562655629
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 11.964182
$ws.Range("H10").Value = 35.892546
$ws.Range("I10").Value = 0.1416244034879904
$ws.Range("J10").Value = 0.1416244034879904
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 47.252316
$ws.Range("N10").Value = 141.756948
$ws.Range("O10").Value = 0.2283395273374399
$ws.Range("P10").Value = 0.2283395273374399
$ws.Range("Q10").Value = 565.335308545512
$ws.Range("R10").Value = 5088.017776909607
$ws.Range("S10").Value = 0.03233844935189462
$ws.Range("T10").Value = 0.03233844935189461
